$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New expense row: "Bag fee on Jet Blue" — $20 "Other" charge from the PR trip,
# dated 2016-01-21 (serial 42390), same as the JetBlue return-flight row above it.
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = 42390
$ws.Range("B23").Value = "Bag fee on Jet Blue"
$ws.Range("I23").Value = 20

# Leave the selection where it was left off after entering the new row.
$ws.Range("E15").Select()
